$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write an exact text value to a cell without Excel
# reinterpreting numeric-looking strings (e.g. "103.00" -> 103).
# A scratch cell is pre-formatted as Text ("@"), given the value,
# then copied with Paste-Values-Only onto the target so the target
# keeps its original (default) cell style/format.
$scratch = $ws.Cells.Item(1, 10)
$scratch.NumberFormat = "@"

function Set-ExactText($cell, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-ExactText $ws.Cells.Item(2, 4) "43.770.61"
Set-ExactText $ws.Cells.Item(2, 5) "  +0.03%  "
Set-ExactText $ws.Cells.Item(3, 4) "2.290.10"
Set-ExactText $ws.Cells.Item(4, 5) "  -0.05%  "
Set-ExactText $ws.Cells.Item(5, 4) "115.38"
Set-ExactText $ws.Cells.Item(5, 5) "  +17.06%  "
Set-ExactText $ws.Cells.Item(6, 4) "268.69"
Set-ExactText $ws.Cells.Item(6, 5) "  -0.78%  "
Set-ExactText $ws.Cells.Item(7, 5) "  +1.15%  "
Set-ExactText $ws.Cells.Item(8, 5) "  +0.12%  "
Set-ExactText $ws.Cells.Item(9, 5) "  +1.57%  "
Set-ExactText $ws.Cells.Item(10, 4) "49.05"
Set-ExactText $ws.Cells.Item(10, 5) "  +8.74%  "
Set-ExactText $ws.Cells.Item(11, 5) "  +0.93%  "
Set-ExactText $ws.Cells.Item(12, 4) "8.97"
Set-ExactText $ws.Cells.Item(12, 5) "  +13.51%  "
Set-ExactText $ws.Cells.Item(13, 4) "0.107"
Set-ExactText $ws.Cells.Item(13, 5) "  +0.34%  "
Set-ExactText $ws.Cells.Item(14, 4) "15.78"
Set-ExactText $ws.Cells.Item(14, 5) "  -0.21%  "
Set-ExactText $ws.Cells.Item(15, 4) "2.634.49"
Set-ExactText $ws.Cells.Item(15, 5) "  -0.24%  "
Set-ExactText $ws.Cells.Item(16, 4) "0.870"
Set-ExactText $ws.Cells.Item(16, 5) "  +1.80%  "
Set-ExactText $ws.Cells.Item(17, 4) "2.290.17"
Set-ExactText $ws.Cells.Item(17, 5) "  -0.25%  "
Set-ExactText $ws.Cells.Item(18, 4) "43.669.74"
Set-ExactText $ws.Cells.Item(18, 5) "  -0.26%  "
Set-ExactText $ws.Cells.Item(19, 5) "  -1.47%  "
Set-ExactText $ws.Cells.Item(20, 5) "  +12.56%  "
Set-ExactText $ws.Cells.Item(21, 4) "72.29"
Set-ExactText $ws.Cells.Item(21, 5) "  +0.05%  "
Set-ExactText $ws.Cells.Item(22, 4) "2.44"
Set-ExactText $ws.Cells.Item(22, 5) "  -0.71%  "
Set-ExactText $ws.Cells.Item(23, 5) "  +0.61%  "
Set-ExactText $ws.Cells.Item(24, 5) "  +8.05%  "
Set-ExactText $ws.Cells.Item(25, 4) "2.93"
Set-ExactText $ws.Cells.Item(25, 5) "  +3.32%  "
Set-ExactText $ws.Cells.Item(26, 4) "11.72"
Set-ExactText $ws.Cells.Item(26, 5) "  +4.00%  "
Set-ExactText $ws.Cells.Item(27, 5) "  -0.04%  "
Set-ExactText $ws.Cells.Item(28, 4) "43.18"
Set-ExactText $ws.Cells.Item(28, 5) "  +13.40%  "
Set-ExactText $ws.Cells.Item(29, 4) "3.93"
Set-ExactText $ws.Cells.Item(29, 5) "  +0.81%  "
Set-ExactText $ws.Cells.Item(31, 5) "  +1.19%  "
Set-ExactText $ws.Cells.Item(32, 4) "173.61"
Set-ExactText $ws.Cells.Item(32, 5) "  -1.75%  "
Set-ExactText $ws.Cells.Item(33, 4) "0.0937"
Set-ExactText $ws.Cells.Item(33, 5) "  +5.21%  "
Set-ExactText $ws.Cells.Item(34, 5) "  -0.96%  "
Set-ExactText $ws.Cells.Item(35, 4) "5.69"
Set-ExactText $ws.Cells.Item(35, 5) "  +4.99%  "
Set-ExactText $ws.Cells.Item(36, 5) "  +0.18%  "
Set-ExactText $ws.Cells.Item(37, 4) "4.81"
Set-ExactText $ws.Cells.Item(37, 5) "  +1.80%  "
Set-ExactText $ws.Cells.Item(38, 5) "  +2.89%  "
Set-ExactText $ws.Cells.Item(39, 5) "  -1.77%  "
Set-ExactText $ws.Cells.Item(40, 4) "3.82"
Set-ExactText $ws.Cells.Item(40, 5) "  +8.48%  "
Set-ExactText $ws.Cells.Item(41, 4) "14.65"
Set-ExactText $ws.Cells.Item(41, 5) "  +20.48%  "
Set-ExactText $ws.Cells.Item(42, 5) "  +15.12%  "
Set-ExactText $ws.Cells.Item(43, 4) "2.41"
Set-ExactText $ws.Cells.Item(43, 5) "  +3.40%  "
Set-ExactText $ws.Cells.Item(44, 4) "0.241"
Set-ExactText $ws.Cells.Item(44, 5) "  +2.09%  "
Set-ExactText $ws.Cells.Item(45, 4) "6.34"
Set-ExactText $ws.Cells.Item(45, 5) "  +21.56%  "
Set-ExactText $ws.Cells.Item(46, 5) "  +0.15%  "
Set-ExactText $ws.Cells.Item(47, 5) "  +2.49%  "
Set-ExactText $ws.Cells.Item(48, 5) "  -0.75%  "
Set-ExactText $ws.Cells.Item(49, 4) "103.00"
Set-ExactText $ws.Cells.Item(49, 5) "  +4.63%  "
Set-ExactText $ws.Cells.Item(50, 5) "  +3.91%  "
Set-ExactText $ws.Cells.Item(51, 4) "0.100"
Set-ExactText $ws.Cells.Item(51, 5) "  -1.83%  "

$scratch.Clear()
$excel.CutCopyMode = $false
Write-Host "Updated cryptos list"
